$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "cm"
$ws.Range("C1").Value = "link_it"
# B1 ("date") and D1 ("topic") are unchanged.

# --- Column A: module codes 001..020 -> cm001..cm020 ---
for ($i = 1; $i -le 20; $i++) {
    $row = $i + 1
    $code = "{0:D3}" -f $i
    $ws.Cells.Item($row, 1).Value = "cm$code"
}

# --- Column D: topic text updates (only rows 2 and 3 change) ---
$ws.Range("D2").Value = "Introduction to computational social science, basic principles of programming, and Python"
$ws.Range("D3").Value = "Loops, conditionals, and functions"

# --- Selection: whole used range instead of single active cell ---
$ws.Range("A1:D21").Select()
